# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel auto-convert a
# numeric-looking string like "233.12" into a Number cell), while leaving
# the cell's style/number-format exactly as it was found.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) '92.083.52'
Set-TextValue $ws.Cells.Item(2, 5) '  +0.65%  '

Set-TextValue $ws.Cells.Item(3, 4) '3.094.58'
Set-TextValue $ws.Cells.Item(3, 5) '  -1.77%  '

Set-TextValue $ws.Cells.Item(4, 5) '  -0.02%  '

Set-TextValue $ws.Cells.Item(5, 4) '233.12'
Set-TextValue $ws.Cells.Item(5, 5) '  -3.35%  '

Set-TextValue $ws.Cells.Item(6, 4) '610.72'
Set-TextValue $ws.Cells.Item(6, 5) '  -1.37%  '

Set-TextValue $ws.Cells.Item(7, 4) '1.08'
Set-TextValue $ws.Cells.Item(7, 5) '  -5.17%  '

Set-TextValue $ws.Cells.Item(8, 4) '0.384'
Set-TextValue $ws.Cells.Item(8, 5) '  +1.90%  '

Set-TextValue $ws.Cells.Item(9, 5) '  -0.05%  '

Set-TextValue $ws.Cells.Item(10, 4) '3.088.29'
Set-TextValue $ws.Cells.Item(10, 5) '  -1.93%  '

Set-TextValue $ws.Cells.Item(11, 4) '0.767'
Set-TextValue $ws.Cells.Item(11, 5) '  +2.47%  '

Set-TextValue $ws.Cells.Item(12, 5) '  -4.11%  '

Set-TextValue $ws.Cells.Item(13, 4) '0.0000242'
Set-TextValue $ws.Cells.Item(13, 5) '  -2.84%  '

Set-TextValue $ws.Cells.Item(14, 4) '91.891.78'
Set-TextValue $ws.Cells.Item(14, 5) '  +0.79%  '

Set-TextValue $ws.Cells.Item(15, 4) '33.52'
Set-TextValue $ws.Cells.Item(15, 5) '  -4.80%  '

Set-TextValue $ws.Cells.Item(16, 4) '5.37'
Set-TextValue $ws.Cells.Item(16, 5) '  -4.33%  '

Set-TextValue $ws.Cells.Item(17, 4) '3.671.77'
Set-TextValue $ws.Cells.Item(17, 5) '  -1.65%  '

Set-TextValue $ws.Cells.Item(18, 4) '3.023.13'
Set-TextValue $ws.Cells.Item(18, 5) '  -4.19%  '

Set-TextValue $ws.Cells.Item(19, 4) '3.78'
Set-TextValue $ws.Cells.Item(19, 5) '  +0.60%  '

Set-TextValue $ws.Cells.Item(20, 4) '14.38'
Set-TextValue $ws.Cells.Item(20, 5) '  -4.77%  '

Set-TextValue $ws.Cells.Item(21, 4) '5.76'
Set-TextValue $ws.Cells.Item(21, 5) '  -4.02%  '

Set-TextValue $ws.Cells.Item(22, 4) '434.42'
Set-TextValue $ws.Cells.Item(22, 5) '  -5.00%  '

Set-TextValue $ws.Cells.Item(23, 4) '9.03'
Set-TextValue $ws.Cells.Item(23, 5) '  -1.34%  '

Set-TextValue $ws.Cells.Item(24, 4) '0.0000196'
Set-TextValue $ws.Cells.Item(24, 5) '  -4.93%  '

Set-TextValue $ws.Cells.Item(25, 4) '5.56'
Set-TextValue $ws.Cells.Item(25, 5) '  -6.42%  '

Set-TextValue $ws.Cells.Item(26, 4) '84.92'
Set-TextValue $ws.Cells.Item(26, 5) '  -4.81%  '

Set-TextValue $ws.Cells.Item(27, 4) '11.27'
Set-TextValue $ws.Cells.Item(27, 5) '  -5.27%  '

Set-TextValue $ws.Cells.Item(28, 4) '3.253.04'
Set-TextValue $ws.Cells.Item(28, 5) '  -2.02%  '

Set-TextValue $ws.Cells.Item(29, 5) '  -0.07%  '

Set-TextValue $ws.Cells.Item(30, 2) 'Cronos'
Set-TextValue $ws.Cells.Item(30, 3) 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Cells.Item(30, 4) '0.176'
Set-TextValue $ws.Cells.Item(30, 5) '  +4.37%  '

Set-TextValue $ws.Cells.Item(31, 2) 'Hedera'
Set-TextValue $ws.Cells.Item(31, 3) 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Cells.Item(31, 4) '0.126'
Set-TextValue $ws.Cells.Item(31, 5) '  -15.83%  '

Set-TextValue $ws.Cells.Item(32, 4) '0.230'
Set-TextValue $ws.Cells.Item(32, 5) '  -3.31%  '

Set-TextValue $ws.Cells.Item(33, 5) '  -36.64%  '

Set-TextValue $ws.Cells.Item(34, 4) '9.08'
Set-TextValue $ws.Cells.Item(34, 5) '  -3.92%  '

Set-TextValue $ws.Cells.Item(35, 4) '7.82'
Set-TextValue $ws.Cells.Item(35, 5) '  +3.89%  '

Set-TextValue $ws.Cells.Item(36, 5) '  -11.99%  '

Set-TextValue $ws.Cells.Item(37, 4) '25.35'
Set-TextValue $ws.Cells.Item(37, 5) '  -4.57%  '

Set-TextValue $ws.Cells.Item(38, 4) '3.88'
Set-TextValue $ws.Cells.Item(38, 5) '  +0.06%  '

Set-TextValue $ws.Cells.Item(39, 4) '1.88'
Set-TextValue $ws.Cells.Item(39, 5) '  -3.53%  '

Set-TextValue $ws.Cells.Item(40, 5) '  +7.57%  '

Set-TextValue $ws.Cells.Item(41, 4) '466.57'
Set-TextValue $ws.Cells.Item(41, 5) '  -5.64%  '

Set-TextValue $ws.Cells.Item(42, 4) '1.27'
Set-TextValue $ws.Cells.Item(42, 5) '  -4.87%  '

Set-TextValue $ws.Cells.Item(43, 4) '0.432'
Set-TextValue $ws.Cells.Item(43, 5) '  -3.90%  '

Set-TextValue $ws.Cells.Item(44, 4) '3.26'
Set-TextValue $ws.Cells.Item(44, 5) '  -4.85%  '

Set-TextValue $ws.Cells.Item(46, 4) '160.14'
Set-TextValue $ws.Cells.Item(46, 5) '  +2.38%  '

Set-TextValue $ws.Cells.Item(47, 4) '0.677'
Set-TextValue $ws.Cells.Item(47, 5) '  -5.51%  '

Set-TextValue $ws.Cells.Item(48, 4) '1.82'
Set-TextValue $ws.Cells.Item(48, 5) '  -6.23%  '

Set-TextValue $ws.Cells.Item(49, 2) 'ImmutableX'
Set-TextValue $ws.Cells.Item(49, 3) 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(49, 4) '1.32'
Set-TextValue $ws.Cells.Item(49, 5) '  -3.41%  '

Set-TextValue $ws.Cells.Item(50, 2) 'OKB'
Set-TextValue $ws.Cells.Item(50, 3) 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Cells.Item(50, 4) '43.77'
Set-TextValue $ws.Cells.Item(50, 5) '  -0.66%  '

Set-TextValue $ws.Cells.Item(51, 2) 'VeChain'
Set-TextValue $ws.Cells.Item(51, 3) 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Cells.Item(51, 4) '0.0323'
Set-TextValue $ws.Cells.Item(51, 5) '  -1.85%  '
